# Updates the cryptos price/volume snapshot on Sheet1.
# Price (D) and Volume(1h) (E) columns are stored as plain text in the
# workbook (not numbers), so each cell is briefly switched to a text
# number format before the assignment (otherwise Excel would parse
# strings like "575.51" as a float) and then restored to the "Normal"
# style so no residual formatting/style change is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '64.104.07'
$r.Style = "Normal"
$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  +0.04%  '
$r.Style = "Normal"
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '2.758.75'
$r.Style = "Normal"
$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  -0.54%  '
$r.Style = "Normal"
$r = $ws.Range('E4')
$r.NumberFormat = "@"
$r.Value = '  +0.06%  '
$r.Style = "Normal"
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '575.51'
$r.Style = "Normal"
$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  -1.90%  '
$r.Style = "Normal"
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '159.43'
$r.Style = "Normal"
$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  -1.22%  '
$r.Style = "Normal"
$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  +0.12%  '
$r.Style = "Normal"
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '0.601'
$r.Style = "Normal"
$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  -2.99%  '
$r.Style = "Normal"
$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  -3.23%  '
$r.Style = "Normal"
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '5.89'
$r.Style = "Normal"
$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  -13.23%  '
$r.Style = "Normal"
$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  +3.37%  '
$r.Style = "Normal"
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '0.386'
$r.Style = "Normal"
$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  -3.23%  '
$r.Style = "Normal"
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '3.250.27'
$r.Style = "Normal"
$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  -0.65%  '
$r.Style = "Normal"
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '27.04'
$r.Style = "Normal"
$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  -2.10%  '
$r.Style = "Normal"
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '63.670.34'
$r.Style = "Normal"
$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  -0.52%  '
$r.Style = "Normal"
$r = $ws.Range('E16')
$r.NumberFormat = "@"
$r.Value = '  -5.60%  '
$r.Style = "Normal"
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '2.767.37'
$r.Style = "Normal"
$r = $ws.Range('E17')
$r.NumberFormat = "@"
$r.Value = '  -0.57%  '
$r.Style = "Normal"
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '12.12'
$r.Style = "Normal"
$r = $ws.Range('E18')
$r.NumberFormat = "@"
$r.Value = '  -2.38%  '
$r.Style = "Normal"
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '4.85'
$r.Style = "Normal"
$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  -3.69%  '
$r.Style = "Normal"
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '359.79'
$r.Style = "Normal"
$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  -1.98%  '
$r.Style = "Normal"
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '6.65'
$r.Style = "Normal"
$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  -6.02%  '
$r.Style = "Normal"
$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  -0.25%  '
$r.Style = "Normal"
$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '0.529'
$r.Style = "Normal"
$r = $ws.Range('E23')
$r.NumberFormat = "@"
$r.Value = '  -7.31%  '
$r.Style = "Normal"
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '65.16'
$r.Style = "Normal"
$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  -3.65%  '
$r.Style = "Normal"
$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  -3.36%  '
$r.Style = "Normal"
$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  -2.49%  '
$r.Style = "Normal"
$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  +0.09%  '
$r.Style = "Normal"
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '0.0₃0905'
$r.Style = "Normal"
$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  -6.89%  '
$r.Style = "Normal"
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '7.40'
$r.Style = "Normal"
$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  +1.63%  '
$r.Style = "Normal"
$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  -3.99%  '
$r.Style = "Normal"
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '1.33'
$r.Style = "Normal"
$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  +4.15%  '
$r.Style = "Normal"
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '170.60'
$r.Style = "Normal"
$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  -1.34%  '
$r.Style = "Normal"
$r = $ws.Range('B33')
$r.NumberFormat = "@"
$r.Value = 'EthereumClassic'
$r.Style = "Normal"
$r = $ws.Range('C33')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$r.Style = "Normal"
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '20.23'
$r.Style = "Normal"
$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  -2.91%  '
$r.Style = "Normal"
$r = $ws.Range('B34')
$r.NumberFormat = "@"
$r.Value = 'NEARProtocol'
$r.Style = "Normal"
$r = $ws.Range('C34')
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$r.Style = "Normal"
$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '4.94'
$r.Style = "Normal"
$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  -3.52%  '
$r.Style = "Normal"
$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  -1.37%  '
$r.Style = "Normal"
$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  +0.08%  '
$r.Style = "Normal"
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '1.81'
$r.Style = "Normal"
$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  -1.36%  '
$r.Style = "Normal"
$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  -2.25%  '
$r.Style = "Normal"
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '351.43'
$r.Style = "Normal"
$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  +2.79%  '
$r.Style = "Normal"
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '6.30'
$r.Style = "Normal"
$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  +0.04%  '
$r.Style = "Normal"
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '4.18'
$r.Style = "Normal"
$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  -2.40%  '
$r.Style = "Normal"
$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '39.10'
$r.Style = "Normal"
$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  -2.64%  '
$r.Style = "Normal"
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '21.61'
$r.Style = "Normal"
$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  -3.89%  '
$r.Style = "Normal"
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '21.98'
$r.Style = "Normal"
$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  -3.12%  '
$r.Style = "Normal"
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '0.0588'
$r.Style = "Normal"
$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  -3.76%  '
$r.Style = "Normal"
$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '137.30'
$r.Style = "Normal"
$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  -1.00%  '
$r.Style = "Normal"
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '0.0254'
$r.Style = "Normal"
$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  -2.81%  '
$r.Style = "Normal"
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '0.629'
$r.Style = "Normal"
$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  -3.33%  '
$r.Style = "Normal"
$r = $ws.Range('E49')
$r.NumberFormat = "@"
$r.Value = '  -1.70%  '
$r.Style = "Normal"
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  +0.08%  '
$r.Style = "Normal"
$r = $ws.Range('E51')
$r.NumberFormat = "@"
$r.Value = '  +0.21%  '
$r.Style = "Normal"
